$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 363.6154
$ws.Range("I55").Value = 182.85715
$ws.Range("J55").Value = 574.5
$ws.Range("K55").Value = 182.85715
$ws.Range("L55").Value = 574.5
$ws.Range("M55").Value = 31.14285000000001
$ws.Range("N55").Value = -1002.5

$ws.Range("H62").Value = 79551864
$ws.Range("I62").Value = 31256542
$ws.Range("K62").Value = 31256542
$ws.Range("M62").Value = -31255918

$ws.Range("H65").Value = 79551864
$ws.Range("I65").Value = 31256542
$ws.Range("K65").Value = 156282710
$ws.Range("M65").Value = -156279590

$ws.Range("H107").Value = 2596.8572
$ws.Range("I107").Value = 2596.8572
$ws.Range("K107").Value = 2596.8572
$ws.Range("M107").Value = -676.8571999999999

$ws.Range("H111").Value = 10267.077
$ws.Range("I111").Value = 4523.3335
$ws.Range("J111").Value = 15190.286
$ws.Range("K111").Value = 13570.0005
$ws.Range("L111").Value = 45570.858
$ws.Range("M111").Value = -10503.0005
$ws.Range("N111").Value = -51704.858

$ws.Range("H113").Value = 7144671.5
$ws.Range("I113").Value = 11112722
$ws.Range("K113").Value = 11112722
$ws.Range("M113").Value = -11109468

$ws.Range("H132").Value = 3474287.2
$ws.Range("I132").Value = 1554.4138
$ws.Range("K132").Value = 4663.2414
$ws.Range("M132").Value = -2133.2414

$ws.Range("H137").Value = 15197293
$ws.Range("I137").Value = 4167580
$ws.Range("J137").Value = 37256720
$ws.Range("K137").Value = 12502740
$ws.Range("L137").Value = 111770160
$ws.Range("M137").Value = -12500190
$ws.Range("N137").Value = -111775260

$ws.Range("H138").Value = 1662.16
$ws.Range("I138").Value = 788.86664
$ws.Range("J138").Value = 2972.1
$ws.Range("K138").Value = 2366.59992
$ws.Range("L138").Value = 8916.299999999999
$ws.Range("M138").Value = 2773.40008
$ws.Range("N138").Value = -19196.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 44446176
$ws.Range("I74").Value = 47619748
$ws.Range("J74").Value = 33338678
$ws.Range("K74").Value = 47619748
$ws.Range("L74").Value = 33338678
$ws.Range("M74").Value = -47618874
$ws.Range("N74").Value = -33340426

$ws.Range("H77").Value = 44446176
$ws.Range("I77").Value = 47619748
$ws.Range("J77").Value = 33338678
$ws.Range("K77").Value = 238098740
$ws.Range("L77").Value = 166693390
$ws.Range("M77").Value = -238094372
$ws.Range("N77").Value = -166702126

$ws.Range("H88").Value = 5589.9
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 5589.9
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 5589.9
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -6401.9

$ws.Range("H91").Value = 5589.9
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 5589.9
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 5589.9
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -8397.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1976
$ws.Range("I86").Value = 1976
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1976
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -853
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 1976
$ws.Range("I89").Value = 1976
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 9880
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4264
$ws.Range("N89").ClearContents()

$ws.Range("H134").Value = 24352160
$ws.Range("I134").Value = 27779112
$ws.Range("J134").Value = 8930871
$ws.Range("K134").Value = 83337336
$ws.Range("L134").Value = 26792613
$ws.Range("M134").Value = -83334801
$ws.Range("N134").Value = -26797683

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2967
$ws.Range("I62").Value = 2463.3333
$ws.Range("J62").Value = 7500
$ws.Range("K62").Value = 2463.3333
$ws.Range("L62").Value = 7500
$ws.Range("M62").Value = -1839.3333
$ws.Range("N62").Value = -8748

$ws.Range("H65").Value = 2967
$ws.Range("I65").Value = 2463.3333
$ws.Range("J65").Value = 7500
$ws.Range("K65").Value = 12316.6665
$ws.Range("L65").Value = 37500
$ws.Range("M65").Value = -9196.666499999999
$ws.Range("N65").Value = -43740

$ws.Range("H132").Value = 2432.8
$ws.Range("I132").Value = 1748.875
$ws.Range("J132").Value = 5168.5
$ws.Range("K132").Value = 5246.625
$ws.Range("L132").Value = 15505.5
$ws.Range("M132").Value = -2716.625
$ws.Range("N132").Value = -20565.5

$ws.Range("H134").Value = 2509427.8
$ws.Range("I134").Value = 9860.833000000001
$ws.Range("J134").Value = 10008128
$ws.Range("K134").Value = 29582.499
$ws.Range("L134").Value = 30024384
$ws.Range("M134").Value = -27047.499
$ws.Range("N134").Value = -30029454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3121328
$ws.Range("I70").Value = 1592288.8
$ws.Range("J70").Value = 7581025.5
$ws.Range("K70").Value = 1592288.8
$ws.Range("L70").Value = 7581025.5
$ws.Range("M70").Value = -1592018.8
$ws.Range("N70").Value = -7581565.5

$ws.Range("H73").Value = 3121328
$ws.Range("I73").Value = 1592288.8
$ws.Range("J73").Value = 7581025.5
$ws.Range("K73").Value = 1592288.8
$ws.Range("L73").Value = 7581025.5
$ws.Range("M73").Value = -1591352.8
$ws.Range("N73").Value = -7582897.5

$ws.Range("H113").Value = 21695
$ws.Range("I113").Value = 4228.75
$ws.Range("J113").Value = 44983.332
$ws.Range("K113").Value = 4228.75
$ws.Range("L113").Value = 44983.332
$ws.Range("M113").Value = -2058.75
$ws.Range("N113").Value = -49323.332

$ws.Range("H126").Value = 34000.668
$ws.Range("I126").Value = 50501
$ws.Range("K126").Value = 151503
$ws.Range("M126").Value = -149033

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7632.2
$ws.Range("I22").Value = 4740
$ws.Range("J22").Value = 10524.4
$ws.Range("K22").Value = 4740
$ws.Range("L22").Value = 10524.4
$ws.Range("M22").Value = -4445
$ws.Range("N22").Value = -11114.4

$ws.Range("H27").Value = 7632.2
$ws.Range("I27").Value = 4740
$ws.Range("J27").Value = 10524.4
$ws.Range("K27").Value = 4740
$ws.Range("L27").Value = 10524.4
$ws.Range("M27").Value = -4633
$ws.Range("N27").Value = -10738.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 372.9091
$ws.Range("I113").Value = 222.4
$ws.Range("J113").Value = 498.33334
$ws.Range("K113").Value = 667.2
$ws.Range("L113").Value = 1495.00002
$ws.Range("M113").Value = 1502.8
$ws.Range("N113").Value = -5835.000019999999

$ws.Range("H132").Value = 980455.25
$ws.Range("I132").Value = 2540.2083
$ws.Range("J132").Value = 4333307
$ws.Range("K132").Value = 7620.624899999999
$ws.Range("L132").Value = 12999921
$ws.Range("M132").Value = -5090.624899999999
$ws.Range("N132").Value = -13004981
